$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-05-28 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-29 Thursday", 2) | Out-Null

# Update the math-problem table cells (20 rows x 5 cols), addressed by position
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "20-4="
$t.Cell(1,2).Range.Text = "36+18="
$t.Cell(1,3).Range.Text = "43-42="
$t.Cell(1,4).Range.Text = "78-2="
$t.Cell(1,5).Range.Text = "51+15="
$t.Cell(2,1).Range.Text = "68+0="
$t.Cell(2,2).Range.Text = "44+37="
$t.Cell(2,3).Range.Text = "36+34="
$t.Cell(2,4).Range.Text = "73-31="
$t.Cell(2,5).Range.Text = "76+5="
$t.Cell(3,1).Range.Text = "64+32="
$t.Cell(3,2).Range.Text = "42-27="
$t.Cell(3,3).Range.Text = "50-11="
$t.Cell(3,4).Range.Text = "97-59="
$t.Cell(3,5).Range.Text = "77-65="
$t.Cell(4,1).Range.Text = "76+6="
$t.Cell(4,2).Range.Text = "97-60="
$t.Cell(4,3).Range.Text = "90-46="
$t.Cell(4,4).Range.Text = "78+6="
$t.Cell(4,5).Range.Text = "55-29="
$t.Cell(5,1).Range.Text = "32+64="
$t.Cell(5,2).Range.Text = "51+7="
$t.Cell(5,3).Range.Text = "98-71="
$t.Cell(5,4).Range.Text = "16+9="
$t.Cell(5,5).Range.Text = "92-54="
$t.Cell(6,1).Range.Text = "98-21="
$t.Cell(6,2).Range.Text = "82+2="
$t.Cell(6,3).Range.Text = "34+29="
$t.Cell(6,4).Range.Text = "28+14="
$t.Cell(6,5).Range.Text = "50+4="
$t.Cell(7,1).Range.Text = "86-26="
$t.Cell(7,2).Range.Text = "89-14="
$t.Cell(7,3).Range.Text = "63+35="
$t.Cell(7,4).Range.Text = "86-57="
$t.Cell(7,5).Range.Text = "6+0="
$t.Cell(8,1).Range.Text = "79-13="
$t.Cell(8,2).Range.Text = "74-34="
$t.Cell(8,3).Range.Text = "55+38="
$t.Cell(8,4).Range.Text = "81-28="
$t.Cell(8,5).Range.Text = "87-43="
$t.Cell(9,1).Range.Text = "45-32="
$t.Cell(9,2).Range.Text = "46-42="
$t.Cell(9,3).Range.Text = "49+39="
$t.Cell(9,4).Range.Text = "64+29="
$t.Cell(9,5).Range.Text = "49+13="
$t.Cell(10,1).Range.Text = "34+64="
$t.Cell(10,2).Range.Text = "91-28="
$t.Cell(10,3).Range.Text = "80-41="
$t.Cell(10,4).Range.Text = "27+68="
$t.Cell(10,5).Range.Text = "28+5="
$t.Cell(11,1).Range.Text = "55-35="
$t.Cell(11,2).Range.Text = "45-13="
$t.Cell(11,3).Range.Text = "96-42="
$t.Cell(11,4).Range.Text = "73+9="
$t.Cell(11,5).Range.Text = "1+67="
$t.Cell(12,1).Range.Text = "35-29="
$t.Cell(12,2).Range.Text = "83-65="
$t.Cell(12,3).Range.Text = "24+44="
$t.Cell(12,4).Range.Text = "55+15="
$t.Cell(12,5).Range.Text = "37+17="
$t.Cell(13,1).Range.Text = "47-12="
$t.Cell(13,2).Range.Text = "39+26="
$t.Cell(13,3).Range.Text = "53+3="
$t.Cell(13,4).Range.Text = "13+76="
$t.Cell(13,5).Range.Text = "32+37="
$t.Cell(14,1).Range.Text = "54+5="
$t.Cell(14,2).Range.Text = "13-3="
$t.Cell(14,3).Range.Text = "43+53="
$t.Cell(14,4).Range.Text = "30+59="
$t.Cell(14,5).Range.Text = "62-37="
$t.Cell(15,1).Range.Text = "60+25="
$t.Cell(15,2).Range.Text = "34+62="
$t.Cell(15,3).Range.Text = "9+33="
$t.Cell(15,4).Range.Text = "49-34="
$t.Cell(15,5).Range.Text = "51+10="
$t.Cell(16,1).Range.Text = "97-25="
$t.Cell(16,2).Range.Text = "7+25="
$t.Cell(16,3).Range.Text = "35-2="
$t.Cell(16,4).Range.Text = "52-47="
$t.Cell(16,5).Range.Text = "77-57="
$t.Cell(17,1).Range.Text = "50-3="
$t.Cell(17,2).Range.Text = "2+8="
$t.Cell(17,3).Range.Text = "79-15="
$t.Cell(17,4).Range.Text = "25+8="
$t.Cell(17,5).Range.Text = "42-41="
$t.Cell(18,1).Range.Text = "27+12="
$t.Cell(18,2).Range.Text = "14+17="
$t.Cell(18,3).Range.Text = "50+30="
$t.Cell(18,4).Range.Text = "39-34="
$t.Cell(18,5).Range.Text = "76-38="
$t.Cell(19,1).Range.Text = "9+62="
$t.Cell(19,2).Range.Text = "24+36="
$t.Cell(19,3).Range.Text = "29-14="
$t.Cell(19,4).Range.Text = "3+75="
$t.Cell(19,5).Range.Text = "67-13="
$t.Cell(20,1).Range.Text = "10+10="
$t.Cell(20,2).Range.Text = "87-45="
$t.Cell(20,3).Range.Text = "70-65="
$t.Cell(20,4).Range.Text = "87-87="
$t.Cell(20,5).Range.Text = "61+28="
